# Scheduled market-data refresh: recompute Leve profit columns (H:N)
# for the affected Leve rows across the per-job Sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 33
$ws.Range("H33").Value = 1524.8518
$ws.Range("I33").Value = 888.9048
$ws.Range("J33").Value = 3750.6667
$ws.Range("K33").Value = 888.9048
$ws.Range("L33").Value = 3750.6667
$ws.Range("M33").Value = -659.9048
$ws.Range("N33").Value = -4208.6667

# Row 69
$ws.Range("H69").Value = 7500
$ws.Range("I69").Value = 50000
$ws.Range("J69").Value = 3636.3635
$ws.Range("K69").Value = 150000
$ws.Range("L69").Value = 10909.0905
$ws.Range("M69").Value = -149126
$ws.Range("N69").Value = -12657.0905

# Row 72
$ws.Range("H72").Value = 7500
$ws.Range("I72").Value = 50000
$ws.Range("J72").Value = 3636.3635
$ws.Range("K72").Value = 450000
$ws.Range("L72").Value = 32727.2715
$ws.Range("M72").Value = -445632
$ws.Range("N72").Value = -41463.2715

# Row 125
$ws.Range("H125").Value = 1358004.2
$ws.Range("I125").Value = 1600.6666
$ws.Range("J125").Value = 1629284.9
$ws.Range("K125").Value = 14405.9994
$ws.Range("L125").Value = 14663564.1
$ws.Range("M125").Value = -11945.9994
$ws.Range("N125").Value = -14668484.1

# Row 129
$ws.Range("H129").Value = 291399.9
$ws.Range("J129").Value = 479300.75
$ws.Range("L129").Value = 1437902.25
$ws.Range("N129").Value = -1447902.25

# Row 138
$ws.Range("H138").Value = 2191.468
$ws.Range("J138").Value = 2988.1428
$ws.Range("L138").Value = 8964.428400000001
$ws.Range("N138").Value = -19244.4284

# Row 141
$ws.Range("H141").Value = 4231.7617
$ws.Range("I141").Value = 2133.182
$ws.Range("J141").Value = 6540.2
$ws.Range("K141").Value = 6399.545999999999
$ws.Range("L141").Value = 19620.6
$ws.Range("M141").Value = -1219.545999999999
$ws.Range("N141").Value = -29980.6

$ws = $wb.Worksheets.Item("ARM")
# Row 45
$ws.Range("H45").Value = 323685.28
$ws.Range("I45").Value = 589142.9399999999
$ws.Range("J45").Value = 1343.8572
$ws.Range("K45").Value = 589142.9399999999
$ws.Range("L45").Value = 1343.8572
$ws.Range("M45").Value = -588765.9399999999
$ws.Range("N45").Value = -2097.8572

# Row 132
$ws.Range("H132").Value = 21268774
$ws.Range("I132").Value = 15051088
$ws.Range("J132").Value = 51113670
$ws.Range("K132").Value = 45153264
$ws.Range("L132").Value = 153341010
$ws.Range("M132").Value = -45150734
$ws.Range("N132").Value = -153346070

$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 57752470
$ws.Range("I20").Value = 33347346
$ws.Range("J20").Value = 91032184
$ws.Range("K20").Value = 33347346
$ws.Range("L20").Value = 91032184
$ws.Range("M20").Value = -33347099
$ws.Range("N20").Value = -91032678

# Row 64
$ws.Range("H64").Value = 9035779
$ws.Range("I64").Value = 1000586
$ws.Range("J64").Value = 11905490
$ws.Range("K64").Value = 1000586
$ws.Range("L64").Value = 11905490
$ws.Range("M64").Value = -1000361
$ws.Range("N64").Value = -11905940

# Row 67
$ws.Range("H67").Value = 9035779
$ws.Range("I67").Value = 1000586
$ws.Range("J67").Value = 11905490
$ws.Range("K67").Value = 1000586
$ws.Range("L67").Value = 11905490
$ws.Range("M67").Value = -999806
$ws.Range("N67").Value = -11907050

$ws = $wb.Worksheets.Item("CRP")
# Row 60
$ws.Range("H60").Value = 8366.666999999999
$ws.Range("I60").Value = 5050
$ws.Range("J60").Value = 15000
$ws.Range("K60").Value = 5050
$ws.Range("L60").Value = 15000
$ws.Range("M60").Value = -4539
$ws.Range("N60").Value = -16022

# Row 120
$ws.Range("H120").Value = 116814.5
$ws.Range("J120").Value = 33333
$ws.Range("L120").Value = 33333
$ws.Range("N120").Value = -40591

$ws = $wb.Worksheets.Item("CUL")
# Row 68
$ws.Range("H68").Value = 2618.016
$ws.Range("I68").Value = 664.13043
$ws.Range("J68").Value = 3770.3076
$ws.Range("K68").Value = 1992.39129
$ws.Range("L68").Value = 11310.9228
$ws.Range("M68").Value = -1181.39129
$ws.Range("N68").Value = -12932.9228

# Row 71
$ws.Range("H71").Value = 2618.016
$ws.Range("I71").Value = 664.13043
$ws.Range("J71").Value = 3770.3076
$ws.Range("K71").Value = 5977.173870000001
$ws.Range("L71").Value = 33932.7684
$ws.Range("M71").Value = -1921.173870000001
$ws.Range("N71").Value = -42044.7684

# Row 103
$ws.Range("H103").Value = 1986.6
$ws.Range("I103").Value = 1266.6666
$ws.Range("J103").Value = 3066.5
$ws.Range("K103").Value = 3799.9998
$ws.Range("L103").Value = 9199.5
$ws.Range("M103").Value = -2920.9998
$ws.Range("N103").Value = -10957.5

# Row 113
$ws.Range("H113").Value = 2319.0933
$ws.Range("I113").Value = 2632.4468
$ws.Range("J113").Value = 1793.1072
$ws.Range("K113").Value = 7897.340400000001
$ws.Range("L113").Value = 5379.321599999999
$ws.Range("M113").Value = -5727.340400000001
$ws.Range("N113").Value = -9719.321599999999

# Row 117
$ws.Range("H117").Value = 1572.2354
$ws.Range("I117").Value = 0
$ws.Range("J117").Value = 1572.2354
$ws.Range("K117").Value = 0
$ws.Range("L117").Value = 4716.706200000001
$ws.Range("N117").Value = -11600.7062
$ws.Range("M117").Value = ""

# Row 131
$ws.Range("H131").Value = 26706.77
$ws.Range("I131").Value = 77343.08
$ws.Range("J131").Value = 1388.6154
$ws.Range("K131").Value = 232029.24
$ws.Range("L131").Value = 4165.8462
$ws.Range("M131").Value = -226989.24
$ws.Range("N131").Value = -14245.8462

$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 4444357
$ws.Range("I70").Value = 2227192
$ws.Range("J70").Value = 11372998
$ws.Range("K70").Value = 2227192
$ws.Range("L70").Value = 11372998
$ws.Range("M70").Value = -2226922
$ws.Range("N70").Value = -11373538

# Row 73
$ws.Range("H73").Value = 4444357
$ws.Range("I73").Value = 2227192
$ws.Range("J73").Value = 11372998
$ws.Range("K73").Value = 2227192
$ws.Range("L73").Value = 11372998
$ws.Range("M73").Value = -2226256
$ws.Range("N73").Value = -11374870

# Row 102
$ws.Range("H102").Value = 3039.068
$ws.Range("I102").Value = 3408.6875
$ws.Range("J102").Value = 2053.4167
$ws.Range("K102").Value = 3408.6875
$ws.Range("L102").Value = 2053.4167
$ws.Range("M102").Value = -1786.6875
$ws.Range("N102").Value = -5297.4167

$ws = $wb.Worksheets.Item("LTW")
# Row 40
$ws.Range("H40").Value = 3475685
$ws.Range("I40").Value = 5053496.5
$ws.Range("J40").Value = 4499.8
$ws.Range("K40").Value = 5053496.5
$ws.Range("L40").Value = 4499.8
$ws.Range("M40").Value = -5053360.5
$ws.Range("N40").Value = -4771.8

# Row 55
$ws.Range("H55").Value = 10869748
$ws.Range("I55").Value = 13889066
$ws.Range("J55").Value = 203.8
$ws.Range("K55").Value = 13889066
$ws.Range("L55").Value = 203.8
$ws.Range("M55").Value = -13888893
$ws.Range("N55").Value = -549.8

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 20677420
$ws.Range("I126").Value = 23109942
$ws.Range("J126").Value = 990
$ws.Range("K126").Value = 69329826
$ws.Range("L126").Value = 2970
$ws.Range("M126").Value = -69327356
$ws.Range("N126").Value = -7910
